$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 0.54
$ws.Range("C5").Value = 0.83
$ws.Range("D5").Value = 0.7
$ws.Range("E5").Value = 0.93
$ws.Range("F5").Value = 0.55
